# Fix inconsistencies in DG diagrams:
#  - Refresh the cached "datetimeFigureOut" date field text (11/1/18 -> 11/6/2018)
#    everywhere it appears: Notes Master, Slide Master, and all Slide Layouts.
#  - Rename the "Model" labels to "Data" on slide 6 (inside nested groups) and
#    slide 7 (top-level rounded rectangle).

$p = $ppt.ActivePresentation

$oldDate = "11/1/18"
$newDate = "11/6/2018"

# --- Notes Master: Date Placeholder ---
$nm = $p.NotesMaster
for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
    $shp = $nm.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Slide Master: Date Placeholder ---
$sm = $p.SlideMaster
for ($i = 1; $i -le $sm.Shapes.Count; $i++) {
    $shp = $sm.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Every Slide Layout: Date Placeholder ---
for ($li = 1; $li -le $sm.CustomLayouts.Count; $li++) {
    $lay = $sm.CustomLayouts.Item($li)
    for ($si = 1; $si -le $lay.Shapes.Count; $si++) {
        $shp = $lay.Shapes.Item($si)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- Slide 6: nested "Model" textbox (inside Group 2 -> flattened GroupItems) ---
$s6 = $p.Slides.Item(6)
$grp = $s6.Shapes.Item(1)
for ($gi = 1; $gi -le $grp.GroupItems.Count; $gi++) {
    $item = $grp.GroupItems.Item($gi)
    if ($item.HasTextFrame) {
        if ($item.TextFrame.TextRange.Text -eq "Model") {
            $item.TextFrame.TextRange.Text = "Data"
        }
    }
}

# --- Slide 7: "Model" rounded-rectangle shape ---
$s7 = $p.Slides.Item(7)
for ($i = 1; $i -le $s7.Shapes.Count; $i++) {
    $shp = $s7.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "Model") {
            $shp.TextFrame.TextRange.Text = "Data"
        }
    }
}
